$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing parameter values (DOC_init and POC_init)
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 0.5

# Add new "#ProductionPeriod" section with growing-season start/end days.
# Values are entered in the same order the shared-string table records
# them (ProdEndDay before ProdStartDay) so the resulting workbook matches
# the author's edit exactly.
$ws.Range("A28").Value = "#ProductionPeriod"

$ws.Range("A30").Value = "ProdEndDay"
$ws.Range("A29").Value = "ProdStartDay"

$ws.Range("C29").Value = "JulianDay"
$ws.Range("C30").Value = "JulianDay"

$ws.Range("B29").Value = 135
$ws.Range("B30").Value = 258

$ws.Range("F13").Select()
